# Generate Report for Handoff
# - Flip the "In Translation" status to "Ready for handoff" everywhere it
#   appears (Overview + per-locale sheets) and refresh the handoff/generate
#   timestamps to reflect the new xliff handoff run.
# - Widen the Status/Latest-Handoff-Datetime columns so the longer
#   "Ready for handoff" label isn't clipped.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-07 14:34:02"
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-07 14:33:47"
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-07 14:34:02"
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
